# "Update fix error log" - swap the crawler configuration row from the
# "techable" site to the "ainow" site on the main (URL) sheet, and move
# the active selection as left by the editing session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2: crawler configuration values -------------------------------
$ws.Range("A2").Value = "ainow"

$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1

$ws.Range("G2").Value = "https://ainow.ai/?s=(keyword)"
$ws.Range("H2").Value = "#main > div > article > a"
$ws.Range("I2").Value = "#main > div > article > a > section > h1"
$ws.Range("J2").Value = "#main > article > header > p > time"
$ws.Range("K2").Value = "null"
$ws.Range("L2").Value = "#main > article > section"

$ws.Range("M2").Value = 2
$ws.Range("N2").Value = "page/"
$ws.Range("O2").Value = "https://ainow.ai"
$ws.Range("P2").Value = "?s="
$ws.Range("Q2").Value = 20
$ws.Range("R2").Value = "null"
$ws.Range("S2").Value = 0

$ws.Range("U2").Value = "#main > nav > ul > li"
$ws.Range("V2").Value = "20件"

# --- Selection, as left after the edit ----------------------------------
$ws.Range("F19").Select() | Out-Null
